$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data values (keep header formatting on A1:H1)
$ws.UsedRange.ClearContents()

# Extend header style (bold, border, centered) from H1 to the new I1 column
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$ws.Range("A1").Value = 'Nombre'
$ws.Range("B1").Value = 'Fuente'
$ws.Range("C1").Value = 'Descripción'
$ws.Range("D1").Value = 'Monto'
$ws.Range("E1").Value = 'Fecha cierre'
$ws.Range("F1").Value = 'Estado'
$ws.Range("G1").Value = 'Área de interés'
$ws.Range("H1").Value = 'Enlace'
$ws.Range("I1").Value = 'País objetivo'

# Data rows (18 new funding opportunities)
# Row 2
$ws.Range("A2").Value = 'FONTAGRO - Convocatoria de Propuestas 2026'
$ws.Range("B2").Value = 'FONTAGRO'
$ws.Range("C2").Value = 'Financiamiento para proyectos innovadores que aumenten sosteniblemente la productividad agrícola en América Latina y el Caribe, especialmente en contexto de cambio climático. Proyectos de hasta 36 meses.'
$ws.Range("D2").Value = 'USD 200,000'
$ws.Range("E2").Value = '''2026-03-30'
$ws.Range("F2").Value = 'Abierto'
$ws.Range("G2").Value = 'Innovación Agrícola'
$ws.Range("H2").Value = 'https://www.fontagro.org/convocatoria-2026/'
$ws.Range("I2").Value = 'América Latina y Caribe'

# Row 3
$ws.Range("A3").Value = 'GAFSP - Eighth Call 2025 for Producer Organizations'
$ws.Range("B3").Value = 'GAFSP'
$ws.Range("C3").Value = 'Programa de $38 millones para fortalecer sistemas alimentarios, mejorar resiliencia climática y empoderar comunidades rurales. Enfocado en organizaciones de productores.'
$ws.Range("D3").Value = 'USD 38,000,000'
$ws.Range("E3").Value = '''2025-11-11'
$ws.Range("F3").Value = 'Abierto'
$ws.Range("G3").Value = 'Seguridad Alimentaria'
$ws.Range("H3").Value = 'https://www.gafspfund.org/call2025'
$ws.Range("I3").Value = 'Global'

# Row 4
$ws.Range("A4").Value = 'EUCaN Facility Nourishing Futures - Protección Social y Nutrición'
$ws.Range("B4").Value = 'EuropeAid'
$ws.Range("C4").Value = 'Financiamiento para sistemas agroalimentarios sostenibles en el Caribe. Enfoque en protección social y nutrición.'
$ws.Range("D4").Value = 'EUR 500,000'
$ws.Range("E4").Value = '''2026-01-20'
$ws.Range("F4").Value = 'Abierto'
$ws.Range("G4").Value = 'Nutrición y Seguridad Alimentaria'
$ws.Range("H4").Value = 'https://www.developmentaid.org/grants/eucan'
$ws.Range("I4").Value = 'Caribe'

# Row 5
$ws.Range("A5").Value = 'Conservation Food and Health Foundation Grants'
$ws.Range("B5").Value = 'Conservation Foundation'
$ws.Range("C5").Value = 'Grants para proyectos en África, Asia, América Latina y Medio Oriente enfocados en producción alimentaria, protección ambiental y salud pública.'
$ws.Range("D5").Value = 'USD 25,000 - 50,000'
$ws.Range("E5").Value = '''2025-06-15'
$ws.Range("F5").Value = 'Abierto'
$ws.Range("G5").Value = 'Producción Alimentaria Sostenible'
$ws.Range("H5").Value = 'https://www.conservationfoodhealth.org/grants'
$ws.Range("I5").Value = 'América Latina'

# Row 6
$ws.Range("A6").Value = 'Rockefeller Foundation - Regenerative Agriculture Initiative'
$ws.Range("B6").Value = 'Rockefeller Foundation'
$ws.Range("C6").Value = 'Iniciativa de $100 millones para construir mercados para producción regenerativa/agroecológica con enfoque en Brasil y América Latina.'
$ws.Range("D6").Value = 'USD 100,000,000'
$ws.Range("E6").Value = '''2026-06-30'
$ws.Range("F6").Value = 'Abierto'
$ws.Range("G6").Value = 'Agricultura Regenerativa'
$ws.Range("H6").Value = 'https://www.rockefellerfoundation.org/regenerative-agriculture'
$ws.Range("I6").Value = 'Brasil y América Latina'

# Row 7
$ws.Range("A7").Value = 'IFAD - Digital Rural Inclusion and Youth Empowerment'
$ws.Range("B7").Value = 'IFAD'
$ws.Range("C7").Value = 'Provisión de tecnologías TIC para inclusión rural digital y empoderamiento juvenil en comunidades rurales.'
$ws.Range("D7").Value = 'USD 500,000'
$ws.Range("E7").Value = '''2025-09-30'
$ws.Range("F7").Value = 'Abierto'
$ws.Range("G7").Value = 'Inclusión Digital Rural'
$ws.Range("H7").Value = 'https://www.ifad.org/tenders'
$ws.Range("I7").Value = 'Global'

# Row 8
$ws.Range("A8").Value = 'FAO Chile - Centro de Semillas Huillilemu'
$ws.Range("B8").Value = 'FAO'
$ws.Range("C8").Value = 'Construcción del Centro de Semillas Huillilemu en la Región de Los Ríos, Chile. Proyecto de infraestructura agrícola.'
$ws.Range("D8").Value = 'USD 800,000'
$ws.Range("E8").Value = '''2025-09-15'
$ws.Range("F8").Value = 'Próximo'
$ws.Range("G8").Value = 'Infraestructura Agrícola'
$ws.Range("H8").Value = 'https://www.fao.org/chile/tenders'
$ws.Range("I8").Value = 'Chile'

# Row 9
$ws.Range("A9").Value = 'UNIDO A2D Facility - Proyectos de Demostración'
$ws.Range("B9").Value = 'UNIDO'
$ws.Range("C9").Value = 'Convocatoria para selección de beneficiarios de grants para implementación de proyectos de demostración A2D en países en desarrollo.'
$ws.Range("D9").Value = 'USD 150,000'
$ws.Range("E9").Value = '''2025-12-31'
$ws.Range("F9").Value = 'Abierto'
$ws.Range("G9").Value = 'Desarrollo Industrial Sostenible'
$ws.Range("H9").Value = 'https://www.unido.org/a2d-facility'
$ws.Range("I9").Value = 'Países en Desarrollo'

# Row 10
$ws.Range("A10").Value = 'AgroLAC 2025 - Productividad Agrícola Sostenible'
$ws.Range("B10").Value = 'BID'
$ws.Range("C10").Value = 'Plataforma multi-donante del BID con The Nature Conservancy para mejorar productividad agrícola y reducir impacto ambiental en América Latina. Fondo total de $50 millones.'
$ws.Range("D10").Value = 'USD 5,000,000'
$ws.Range("E10").Value = '''2025-12-15'
$ws.Range("F10").Value = 'Abierto'
$ws.Range("G10").Value = 'Productividad Agrícola'
$ws.Range("H10").Value = 'https://www.iadb.org/agrolac'
$ws.Range("I10").Value = 'América Latina'

# Row 11
$ws.Range("A11").Value = 'Climate-Smart Agriculture Fund (CSAF)'
$ws.Range("B11").Value = 'NDF/BID'
$ws.Range("C11").Value = 'Fondo de financiamiento concesional para atraer inversión del sector privado hacia agricultura sostenible, silvicultura y desarrollo de pastizales en la región.'
$ws.Range("D11").Value = 'EUR 5,000,000'
$ws.Range("E11").Value = '''2026-03-15'
$ws.Range("F11").Value = 'Abierto'
$ws.Range("G11").Value = 'Agricultura Climáticamente Inteligente'
$ws.Range("H11").Value = 'https://www.ndf.int/csaf'
$ws.Range("I11").Value = 'América Latina y Caribe'

# Row 12
$ws.Range("A12").Value = 'América Latina y el Caribe Sin Hambre 2025'
$ws.Range("B12").Value = 'FAO'
$ws.Range("C12").Value = 'Programa de cooperación Brasil-FAO enfocado en seguridad alimentaria y nutricional, reducción de pobreza y asistencia técnica. Énfasis en agricultura familiar.'
$ws.Range("D12").Value = 'USD 2,000,000'
$ws.Range("E12").Value = '''2025-12-31'
$ws.Range("F12").Value = 'Abierto'
$ws.Range("G12").Value = 'Seguridad Alimentaria'
$ws.Range("H12").Value = 'https://www.fao.org/alc-sin-hambre'
$ws.Range("I12").Value = 'América Latina y Caribe'

# Row 13
$ws.Range("A13").Value = 'Gates Foundation - Agricultural Development Grant'
$ws.Range("B13").Value = 'Gates Foundation'
$ws.Range("C13").Value = 'Financiamiento para proyectos de desarrollo agrícola con enfoque en pequeños productores y sistemas alimentarios sostenibles.'
$ws.Range("D13").Value = 'USD 1,000,000'
$ws.Range("E13").Value = '''2025-03-25'
$ws.Range("F13").Value = 'Próximo'
$ws.Range("G13").Value = 'Desarrollo Agrícola'
$ws.Range("H13").Value = 'https://www.developmentaid.org/grants/gates'
$ws.Range("I13").Value = 'Global'

# Row 14
$ws.Range("A14").Value = 'UN Global Indigenous Youth Forum 2026 - Sistemas Alimentarios'
$ws.Range("B14").Value = 'ONU'
$ws.Range("C14").Value = 'Convocatoria para jóvenes indígenas enfocada en sistemas alimentarios y conocimiento tradicional, biodiversidad, restauración de ecosistemas y resiliencia climática.'
$ws.Range("D14").Value = 'USD 50,000'
$ws.Range("E14").Value = '''2026-06-30'
$ws.Range("F14").Value = 'Abierto'
$ws.Range("G14").Value = 'Conocimiento Indígena'
$ws.Range("H14").Value = 'https://www.fao.org/ungiyf'
$ws.Range("I14").Value = 'Global'

# Row 15
$ws.Range("A15").Value = 'CELAC Plan SAN 2024-2030 - Erradicación del Hambre'
$ws.Range("B15").Value = 'CELAC'
$ws.Range("C15").Value = 'Marco regional para alcanzar los ODS relacionados con el hambre y la malnutrición. Financiamiento para proyectos nacionales alineados.'
$ws.Range("D15").Value = 'USD 500,000'
$ws.Range("E15").Value = '''2026-12-31'
$ws.Range("F15").Value = 'Abierto'
$ws.Range("G15").Value = 'Erradicación del Hambre'
$ws.Range("H15").Value = 'https://www.cepal.org/celac-san'
$ws.Range("I15").Value = 'América Latina y Caribe'

# Row 16
$ws.Range("A16").Value = 'FIA - Convocatoria Nacional de Innovación 2026'
$ws.Range("B16").Value = 'FIA'
$ws.Range("C16").Value = 'Financiamiento para proyectos de innovación agrícola en Chile. Enfoque en digitalización, sustentabilidad y adaptación al cambio climático.'
$ws.Range("D16").Value = 'CLP 200,000,000'
$ws.Range("E16").Value = '''2026-04-30'
$ws.Range("F16").Value = 'Abierto'
$ws.Range("G16").Value = 'Innovación Agrícola'
$ws.Range("H16").Value = 'https://www.fia.cl/convocatorias'
$ws.Range("I16").Value = 'Chile'

# Row 17
$ws.Range("A17").Value = 'INDAP - Programa de Desarrollo Rural 2026'
$ws.Range("B17").Value = 'INDAP'
$ws.Range("C17").Value = 'Programa para pequeños agricultores chilenos. Incluye asistencia técnica, financiamiento y capacitación.'
$ws.Range("D17").Value = 'CLP 50,000,000'
$ws.Range("E17").Value = '''2026-03-31'
$ws.Range("F17").Value = 'Abierto'
$ws.Range("G17").Value = 'Desarrollo Rural'
$ws.Range("H17").Value = 'https://www.indap.cl/programas'
$ws.Range("I17").Value = 'Chile'

# Row 18
$ws.Range("A18").Value = 'CORFO Innova - Agroindustria Sustentable'
$ws.Range("B18").Value = 'CORFO'
$ws.Range("C18").Value = 'Financiamiento para proyectos de innovación en agroindustria chilena con enfoque en sustentabilidad y economía circular.'
$ws.Range("D18").Value = 'CLP 300,000,000'
$ws.Range("E18").Value = '''2026-05-15'
$ws.Range("F18").Value = 'Abierto'
$ws.Range("G18").Value = 'Agroindustria'
$ws.Range("H18").Value = 'https://www.corfo.cl/innova-agro'
$ws.Range("I18").Value = 'Chile'

# Row 19
$ws.Range("A19").Value = 'ANID FONDECYT - Investigación Agrícola 2026'
$ws.Range("B19").Value = 'ANID'
$ws.Range("C19").Value = 'Financiamiento para investigación científica en áreas agrícolas, incluyendo biotecnología, recursos hídricos y cambio climático.'
$ws.Range("D19").Value = 'CLP 150,000,000'
$ws.Range("E19").Value = '''2026-06-30'
$ws.Range("F19").Value = 'Abierto'
$ws.Range("G19").Value = 'Investigación Agrícola'
$ws.Range("H19").Value = 'https://www.anid.cl/fondecyt'
$ws.Range("I19").Value = 'Chile'

# The leading apostrophe above forces 'Fecha cierre' to stay plain text instead of
# being auto-parsed into a date serial; ClearFormats() then drops the resulting
# quotePrefix style so the cells end up unstyled, matching the source data.
$ws.Range("E2:E19").ClearFormats()
